$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new Density column (D1): bold font, centered, bordered - like other headers
$ws.Range("D1").Value = "Density"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").Borders.LineStyle = 1
$ws.Range("D1").HorizontalAlignment = -4108

# Update existing C column (avg_weight) cells to add center alignment (border already present)
$ws.Range("C2:C46").HorizontalAlignment = -4108

# Apply border + center alignment to the new D column data range
$ws.Range("D2:D46").Borders.LineStyle = 1
$ws.Range("D2:D46").HorizontalAlignment = -4108

# Density values per foodstuff row, in row order 2..46
$densities = @("0.92","0.70","0.96","1.00","0.73","1.15","0.95","1.00","0.65","0.45","0.36","1.12","0.55","0.65","0.63","0.58","1.00","0.65","0.85","0.98","1.09","0.85","0.98","0.98","0.23","0.98","0.97","0.65","0.94","0.51","0.98","0.95","0.99","0.85","0.99","0.99","1.05","0.68","0.95","0.95","0.98","0.99","0.60","0.98","0.65")

for ($i = 0; $i -lt $densities.Length; $i++) {
    $row = $i + 2
    $addr = "D" + $row
    $ws.Range($addr).Value = "'" + $densities[$i]
}

# Update selection to match the authored workbook (no frozen topLeftCell, active cell C6)
$ws.Range("C6").Select()
